# Patient details.xlsx - "People added up to 80"
# Fills in First Name (column C) and Last Name (column D) for rows 47-80,
# which previously only had NHS Number / Title / Address / phone / blood
# group / allergy data. Names are mostly derived from The Wolf of Wall
# Street cast. Cells are written in the same order the author originally
# typed them in (not strictly row-by-row) so that newly-introduced shared
# strings line up with the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C47").Value = 'Leonardo'
$ws.Range("D47").Value = 'DiCaprio'

$ws.Range("C48").Value = 'Jonah'
$ws.Range("D48").Value = 'Hill'

$ws.Range("C49").Value = 'Matthew'
$ws.Range("D49").Value = 'McConaughey'

$ws.Range("C56").Value = 'Kyle'
$ws.Range("D56").Value = 'Chandler'

$ws.Range("C57").Value = 'Rob'
$ws.Range("D57").Value = 'Reiner'

$ws.Range("C59").Value = 'Jon'
$ws.Range("D59").Value = 'Bernthal'

$ws.Range("C61").Value = 'Jon'
$ws.Range("D61").Value = 'Favreau'

$ws.Range("C62").Value = 'Jean'
$ws.Range("D62").Value = 'Dujardin'

$ws.Range("C60").Value = 'Joanna'
$ws.Range("D60").Value = 'Lumley'

$ws.Range("C58").Value = 'Cristin'
$ws.Range("D58").Value = 'Milioti'

$ws.Range("C55").Value = 'Christine'
$ws.Range("D55").Value = 'Ebersole'

$ws.Range("C66").Value = 'Shea'
$ws.Range("D66").Value = 'Whigham'

$ws.Range("C50").Value = 'Katarina'
$ws.Range("D50").Value = 'Cas'

$ws.Range("C68").Value = 'Kenneth'
$ws.Range("D68").Value = 'Choi'

$ws.Range("C69").Value = 'Brian'
$ws.Range("D69").Value = 'Sacca'

$ws.Range("C70").Value = 'Henry'
$ws.Range("D70").Value = 'Zebrowski'

$ws.Range("C73").Value = 'Ethan'
$ws.Range("D73").Value = 'Suplee'

$ws.Range("C74").Value = 'Barry'
$ws.Range("D74").Value = 'Rothbart'

$ws.Range("C63").Value = 'Megan'
$ws.Range("D63").Value = 'McKenzie'

$ws.Range("C77").Value = 'Jake'
$ws.Range("D77").Value = 'Hoffman'

$ws.Range("C64").Value = 'Amy'
$ws.Range("D64").Value = 'Cash'

$ws.Range("C65").Value = 'Stephanie'
$ws.Range("D65").Value = 'McCoy'

$ws.Range("C71").Value = 'Ashley'
$ws.Range("D71").Value = 'Atkinson'

$ws.Range("C80").Value = 'Ted'
$ws.Range("D80").Value = 'Griffin'

$ws.Range("C79").Value = 'Edward'
$ws.Range("D79").Value = 'Hermann'

$ws.Range("C78").Value = 'Stephen'
$ws.Range("D78").Value = 'O''Neill'

$ws.Range("C76").Value = 'Natasha'
$ws.Range("D76").Value = 'Thomas'

$ws.Range("C75").Value = 'Sandra'
$ws.Range("D75").Value = 'Nelson'

$ws.Range("C72").Value = 'Carla'
$ws.Range("D72").Value = 'Flaherty'

$ws.Range("C67").Value = 'Madison'
$ws.Range("D67").Value = 'McKinley'

$ws.Range("C51").Value = 'Diedre'
$ws.Range("D51").Value = 'Reimond'

$ws.Range("C52").Value = 'Kelly'
$ws.Range("D52").Value = 'Malloy'

$ws.Range("C53").Value = 'Stephanie'
$ws.Range("D53").Value = 'Witting'

$ws.Range("C54").Value = 'Tess'
$ws.Range("D54").Value = 'Gillis'

# Match the author's final cursor / scroll position from the diff
# (topLeftCell moved to A48, active selection to D54).
$ws.Range("A48").Select() | Out-Null
$ws.Range("D54").Select() | Out-Null
